# Scrum Master: Files update
# Update the Sprint 1 burndown-chart tracker: shift the sprint dates forward,
# log a bit of work done on Day 4 (tasks #1 and #2), bump the estimate for
# task #2, and drop task #9 (it never got logged / was removed from scope).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: shift the daily dates forward (Day 0 .. Day 6) ---
$ws.Range("E4").Value = 45220
$ws.Range("F4").Value = 45221
$ws.Range("G4").Value = 45222
$ws.Range("H4").Value = 45223
$ws.Range("I4").Value = 45224
$ws.Range("J4").Value = 45225
$ws.Range("K4").Value = 45226

# --- Task #2 (row 7): initial estimate grows from 2 to 3 ---
$ws.Range("D7").Value = 3

# --- Log effort on Day 4 (column H) for tasks #1 and #2 ---
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1

# --- Task #9 (row 14) is removed: clear the description and estimate ---
$ws.Range("C14").ClearContents()
$ws.Range("D14").Value = ""

# --- Update the active selection left in the sheet view ---
$ws.Range("G9").Select()
